# V2.0 Upgrade Checklist - "Started adding async methods"
#
# Adds five new tracking columns (C:G) to the "Methods" sheet, each marked
# "Done" (the existing green "Good" style) for the Entity/Answer-group rows
# that were already complete, plus two brand-new "Done" rows (Synchronous
# wrapper rows) for the Answers entity and final catch-all entity.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Methods" sheet
$ws.Activate()

# --- New header row (columns C:G) ---------------------------------------
$ws.Range("C1").Value = "Synchronus"
$ws.Range("D1").Value = "Asynchronus"
$ws.Range("E1").Value = ".Net 3.5"
$ws.Range("F1").Value = "Integration Test"
$ws.Range("G1").Value = "Documentation"

# --- Rows that get a brand-new "Done" marker in column A (Entity rows) --
$newEntityDoneRows = @(2, 3, 88)
foreach ($r in $newEntityDoneRows) {
    $cell = $ws.Range("A$r")
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# --- Rows that get "Done" added in the new Synchronus/Asynchronus cols --
$syncAsyncDoneRows = @(2, 3, 4, 6, 7, 8, 9, 10, 11, 13, 14, 81, 82, 83, 84)
foreach ($r in $syncAsyncDoneRows) {
    $ws.Range("C$r").Value = "Done"
    $ws.Range("D$r").Value = "Done"
}

# --- Size the new columns to fit their (now-populated) contents ----------
$ws.Columns("C:G").AutoFit()

# --- Restore the saved selection -----------------------------------------
$ws.Range("D81").Select()
